$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $val)
    if ($val -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        $cell.NumberFormat = "@"
    }
    $cell.Value2 = $val
}

# Rows 2-33: update Price (D) and Volume (E) columns only
Set-TextValue $ws.Range("D2") "26.504.53"
Set-TextValue $ws.Range("E2") "  +4.01%  "
Set-TextValue $ws.Range("D3") "1.734.84"
Set-TextValue $ws.Range("E3") "  +4.33%  "
Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  +0.15%  "
Set-TextValue $ws.Range("D5") "243.87"
Set-TextValue $ws.Range("E5") "  +3.53%  "
Set-TextValue $ws.Range("E6") "  +0.07%  "
Set-TextValue $ws.Range("D7") "0.4790"
Set-TextValue $ws.Range("E7") "  +3.59%  "
Set-TextValue $ws.Range("D8") "0.2662"
Set-TextValue $ws.Range("E8") "  +3.53%  "
Set-TextValue $ws.Range("D9") "0.06223"
Set-TextValue $ws.Range("E9") "  +1.41%  "
Set-TextValue $ws.Range("D10") "1.736.72"
Set-TextValue $ws.Range("E10") "  +4.40%  "
Set-TextValue $ws.Range("D11") "0.07127"
Set-TextValue $ws.Range("E11") "  +2.66%  "
Set-TextValue $ws.Range("D12") "15.71"
Set-TextValue $ws.Range("E12") "  +7.45%  "
Set-TextValue $ws.Range("D13") "0.6134"
Set-TextValue $ws.Range("E13") "  +7.07%  "
Set-TextValue $ws.Range("D14") "4.526"
Set-TextValue $ws.Range("E14") "  +4.61%  "
Set-TextValue $ws.Range("D15") "76.78"
Set-TextValue $ws.Range("E15") "  +2.31%  "
Set-TextValue $ws.Range("E16") "  +0.10%  "
Set-TextValue $ws.Range("D17") "26.509.51"
Set-TextValue $ws.Range("E17") "  +4.02%  "
Set-TextValue $ws.Range("D18") "1.001"
Set-TextValue $ws.Range("E18") "  +0.11%  "
Set-TextValue $ws.Range("D19") "0.000006892"
Set-TextValue $ws.Range("E19") "  +2.85%  "
Set-TextValue $ws.Range("D20") "11.71"
Set-TextValue $ws.Range("E20") "  +3.08%  "
Set-TextValue $ws.Range("D21") "1.959.31"
Set-TextValue $ws.Range("E21") "  +4.50%  "
Set-TextValue $ws.Range("D22") "4.557"
Set-TextValue $ws.Range("E22") "  +3.40%  "
Set-TextValue $ws.Range("D23") "8.867"
Set-TextValue $ws.Range("E23") "  +2.60%  "
Set-TextValue $ws.Range("D24") "5.329"
Set-TextValue $ws.Range("E24") "  +2.28%  "
Set-TextValue $ws.Range("D25") "135.82"
Set-TextValue $ws.Range("E25") "  +0.82%  "
Set-TextValue $ws.Range("D26") "15.33"
Set-TextValue $ws.Range("E26") "  +3.11%  "
Set-TextValue $ws.Range("D27") "1.797"
Set-TextValue $ws.Range("E27") "  +5.23%  "
Set-TextValue $ws.Range("D28") "1.397"
Set-TextValue $ws.Range("E28") "  +2.33%  "
Set-TextValue $ws.Range("D29") "106.42"
Set-TextValue $ws.Range("E29") "  +2.74%  "
Set-TextValue $ws.Range("D30") "3.971"
Set-TextValue $ws.Range("E30") "  +0.59%  "
Set-TextValue $ws.Range("D31") "3.704"
Set-TextValue $ws.Range("E31") "  +3.16%  "
Set-TextValue $ws.Range("D32") "0.07878"
Set-TextValue $ws.Range("E32") "  +2.28%  "
Set-TextValue $ws.Range("D33") "0.04555"
Set-TextValue $ws.Range("E33") "  +5.22%  "

# Row 34: new Frax row inserted; rows below shift down by one
Set-TextValue $ws.Range("B34") "Frax"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D34") "1.000"
Set-TextValue $ws.Range("E34") "  +0.17%  "

# Rows 35-51: shifted content (old row-1 data) with updated Price/Volume
Set-TextValue $ws.Range("B35") "HuobiToken"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D35") "2.618"
Set-TextValue $ws.Range("E35") "  +0.19%  "
Set-TextValue $ws.Range("B36") "ImmutableX"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "0.6342"
Set-TextValue $ws.Range("E36") "  +5.74%  "
Set-TextValue $ws.Range("B37") "ARBITRUM"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D37") "0.9915"
Set-TextValue $ws.Range("E37") "  +5.70%  "
Set-TextValue $ws.Range("B38") "TrustWalletToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "0.9347"
Set-TextValue $ws.Range("E38") "  +2.93%  "
Set-TextValue $ws.Range("B39") "Quant"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D39") "110.89"
Set-TextValue $ws.Range("E39") "  +3.16%  "
Set-TextValue $ws.Range("B40") "MXToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.431"
Set-TextValue $ws.Range("E40") "  +0.53%  "
Set-TextValue $ws.Range("B41") "RenderToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D41") "1.979"
Set-TextValue $ws.Range("E41") "  +8.87%  "
Set-TextValue $ws.Range("B42") "PaxDollar"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D42") "1.005"
Set-TextValue $ws.Range("E42") "  +0.60%  "
Set-TextValue $ws.Range("B43") "VeChain"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D43") "0.01510"
Set-TextValue $ws.Range("E43") "  +3.68%  "
Set-TextValue $ws.Range("B44") "FraxShare"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D44") "5.698"
Set-TextValue $ws.Range("E44") "  +13.98%  "
Set-TextValue $ws.Range("B45") "TheSandbox"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D45") "0.3900"
Set-TextValue $ws.Range("E45") "  +5.23%  "
Set-TextValue $ws.Range("B46") "Aptos"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D46") "6.903"
Set-TextValue $ws.Range("E46") "  +12.97%  "
Set-TextValue $ws.Range("B47") "Algorand"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D47") "0.1192"
Set-TextValue $ws.Range("E47") "  +7.82%  "
Set-TextValue $ws.Range("B48") "Cronos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.05334"
Set-TextValue $ws.Range("E48") "  +1.48%  "
Set-TextValue $ws.Range("B49") "EnergySwap"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "7.869"
Set-TextValue $ws.Range("E49") "  +3.37%  "
Set-TextValue $ws.Range("B50") "Elrond"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D50") "30.77"
Set-TextValue $ws.Range("E50") "  +0.91%  "
Set-TextValue $ws.Range("B51") "NEARProtocol"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.255"
Set-TextValue $ws.Range("E51") "  +5.48%  "

Write-Host "Applied cryptos update"